$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.166.34"
$ws.Range("E2").Value = "  +6.31%  "
$ws.Range("D3").Value = "3.549.89"
$ws.Range("E3").Value = "  +3.57%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "418.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.660"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.03%  "
$ws.Range("D8").Value = "3.539.27"
$ws.Range("E8").Value = "  +3.47%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.782"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.169"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +20.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000292"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +34.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "43.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("E14").Value = "  +8.63%  "
$ws.Range("D15").Value = "4.123.71"
$ws.Range("E15").Value = "  +4.00%  "
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "3.565.72"
$ws.Range("E18").Value = "  +3.81%  "
$ws.Range("B19").Value = "Polygon"
$ws.Range("C19").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.79%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "66.094.49"
$ws.Range("E21").Value = "  +6.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "448.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "90.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("E30").Value = "  +5.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.25%  "
$ws.Range("E32").Value = "  +5.67%  "
$ws.Range("E33").Value = "  -4.53%  "
$ws.Range("E34").Value = "  -2.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0507"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.30%  "
$ws.Range("D39").Value = "0.0₃0736"
$ws.Range("E39").Value = "  +40.61%  "
$ws.Range("E40").Value = "  +10.83%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "147.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("E47").Value = "  -4.66%  "
$ws.Range("E48").Value = "  -3.73%  "
$ws.Range("E49").Value = "  -4.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.148"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "15.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.95%  "
